$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, $text)
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value2 = $text
    $r.ClearFormats()
}

Set-CellText 'D2' '24.958.45'
Set-CellText 'D3' '1.706.52'
Set-CellText 'E3' '  +1.39%  '
Set-CellText 'D4' '1.001'
Set-CellText 'E4' '  +0.06%  '
Set-CellText 'D5' '316.36'
Set-CellText 'E5' '  +0.05%  '
Set-CellText 'D6' '1.000'
Set-CellText 'E6' '  +0.05%  '
Set-CellText 'D7' '0.3953'
Set-CellText 'E7' '  +1.89%  '
Set-CellText 'D8' '0.4031'
Set-CellText 'E8' '  +0.84%  '
Set-CellText 'E9' '  +0.49%  '
Set-CellText 'D10' '52.64'
Set-CellText 'E10' '  +1.11%  '
Set-CellText 'D11' '1.000'
Set-CellText 'E11' '  +0.10%  '
Set-CellText 'D12' '0.08825'
Set-CellText 'E12' '  +1.03%  '
Set-CellText 'D13' '26.15'
Set-CellText 'E13' '  +0.94%  '
Set-CellText 'D14' '7.478'
Set-CellText 'E14' '  +0.20%  '
Set-CellText 'D15' '0.00001356'
Set-CellText 'E15' '  +1.23%  '
Set-CellText 'E16' '  +0.32%  '
Set-CellText 'D17' '1.718.66'
Set-CellText 'E17' '  +2.85%  '
Set-CellText 'D18' '96.25'
Set-CellText 'E18' '  -1.34%  '
Set-CellText 'D19' '0.07182'
Set-CellText 'E19' '  -0.11%  '
Set-CellText 'D20' '20.59'
Set-CellText 'E20' '  +4.90%  '
Set-CellText 'D21' '7.352'
Set-CellText 'E21' '  +1.65%  '
Set-CellText 'D22' '0.9995'
Set-CellText 'E22' '  +0.00%  '
Set-CellText 'D23' '14.47'
Set-CellText 'E23' '  +2.43%  '
Set-CellText 'D24' '24.988.30'
Set-CellText 'E24' '  +2.45%  '
Set-CellText 'D25' '2.979'
Set-CellText 'E25' '  -0.91%  '
Set-CellText 'D26' '2.348'
Set-CellText 'E26' '  +0.34%  '
Set-CellText 'D27' '23.69'
Set-CellText 'E27' '  +5.53%  '
Set-CellText 'D28' '6.244'
Set-CellText 'E28' '  +16.62%  '
Set-CellText 'D29' '161.64'
Set-CellText 'E29' '  -3.51%  '
Set-CellText 'D30' '150.48'
Set-CellText 'E30' '  +9.20%  '
Set-CellText 'D31' '8.437'
Set-CellText 'E31' '  -1.93%  '
Set-CellText 'D32' '2.590'
Set-CellText 'E32' '  +31.44%  '
Set-CellText 'D33' '1.914.50'
Set-CellText 'E33' '  +3.25%  '
Set-CellText 'D34' '0.08577'
Set-CellText 'E34' '  -1.81%  '
Set-CellText 'E35' '  +0.57%  '
Set-CellText 'B36' 'InternetComputer(DFINITY)'
Set-CellText 'C36' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText 'D36' '7.204'
Set-CellText 'E36' '  -1.73%  '
Set-CellText 'B37' 'VeChain'
Set-CellText 'C37' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText 'D37' '0.03116'
Set-CellText 'E37' '  +4.90%  '
Set-CellText 'D38' '0.2861'
Set-CellText 'E38' '  +4.40%  '
Set-CellText 'B39' 'FraxShare'
Set-CellText 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText 'D39' '10.86'
Set-CellText 'E39' '  +1.06%  '
Set-CellText 'B40' 'Stellar'
Set-CellText 'C40' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText 'D40' '0.09531'
Set-CellText 'E40' '  +4.41%  '
Set-CellText 'D41' '0.8290'
Set-CellText 'E41' '  +4.28%  '
Set-CellText 'E42' '  +0.11%  '
Set-CellText 'D43' '1.486'
Set-CellText 'E43' '  +1.22%  '
Set-CellText 'D44' '17.41'
Set-CellText 'E44' '  +0.85%  '
Set-CellText 'D45' '2.695'
Set-CellText 'E45' '  +4.40%  '
Set-CellText 'D46' '0.7404'
Set-CellText 'E46' '  +3.18%  '
Set-CellText 'D47' '4.251'
Set-CellText 'E47' '  -0.16%  '
Set-CellText 'D48' '1.402'
Set-CellText 'E48' '  +0.63%  '
Set-CellText 'D49' '0.08765'
Set-CellText 'E49' '  +9.07%  '
Set-CellText 'E50' '  +0.16%  '
Set-CellText 'D51' '139.22'
Set-CellText 'E51' '  -0.01%  '
